# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the market-price columns (H:N) across
# several worksheets, matching the upstream scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Value updates
$ws.Range("H43").Value = 1943.7142
$ws.Range("I43").Value = 1779
$ws.Range("K43").Value = 1779
$ws.Range("M43").Value = -1710
$ws.Range("H51").Value = 250004500
$ws.Range("I51").Value = 500000000
$ws.Range("J51").Value = 9000
$ws.Range("K51").Value = 500000000
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = -499999516
$ws.Range("N51").Value = -9968
$ws.Range("H57").Value = 64204.5
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 64204.5
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 192613.5
$ws.Range("N57").Value = -193611.5
$ws.Range("H64").Value = 200004320
$ws.Range("J64").Value = 500004300
$ws.Range("L64").Value = 500004300
$ws.Range("N64").Value = -500004796
$ws.Range("H67").Value = 200004320
$ws.Range("J67").Value = 500004300
$ws.Range("L67").Value = 500004300
$ws.Range("N67").Value = -500006016
$ws.Range("H74").Value = 16867.938
$ws.Range("I74").Value = 18529.846
$ws.Range("J74").Value = 9666.333000000001
$ws.Range("K74").Value = 18529.846
$ws.Range("L74").Value = 9666.333000000001
$ws.Range("M74").Value = -17593.846
$ws.Range("N74").Value = -11538.333
$ws.Range("H77").Value = 16867.938
$ws.Range("I77").Value = 18529.846
$ws.Range("J77").Value = 9666.333000000001
$ws.Range("K77").Value = 92649.23000000001
$ws.Range("L77").Value = 48331.665
$ws.Range("M77").Value = -87969.23000000001
$ws.Range("N77").Value = -57691.665
$ws.Range("H138").Value = 3279.1428
$ws.Range("I138").Value = 2684.2222
$ws.Range("K138").Value = 8052.6666
$ws.Range("H141").Value = 1913.5
$ws.Range("I141").Value = 1980.9412
$ws.Range("J141").Value = 767
$ws.Range("K141").Value = 5942.8236
$ws.Range("L141").Value = 2301
$ws.Range("M141").Value = -762.8235999999997

# New cells (previously empty)
$ws.Range("M138").Value = -2912.6666
$ws.Range("N141").Value = -12661

# Cells cleared (now empty)
$ws.Range("M57").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Value updates
$ws.Range("H32").Value = 2763.4
$ws.Range("I32").Value = 2841.4482
$ws.Range("K32").Value = 2841.4482
$ws.Range("M32").Value = -2554.4482
$ws.Range("H61").Value = 2975.9583
$ws.Range("I61").Value = 2089.0625
$ws.Range("K61").Value = 2089.0625
$ws.Range("M61").Value = -1877.0625
$ws.Range("H74").Value = 150526.42
$ws.Range("I74").Value = 177737.5
$ws.Range("J74").Value = 5400.6665
$ws.Range("K74").Value = 177737.5
$ws.Range("L74").Value = 5400.6665
$ws.Range("M74").Value = -176863.5
$ws.Range("N74").Value = -7148.6665
$ws.Range("H77").Value = 150526.42
$ws.Range("I77").Value = 177737.5
$ws.Range("J77").Value = 5400.6665
$ws.Range("K77").Value = 888687.5
$ws.Range("L77").Value = 27003.3325
$ws.Range("M77").Value = -884319.5
$ws.Range("N77").Value = -35739.3325
$ws.Range("H132").Value = 1919.0227
$ws.Range("I132").Value = 1667.9333
$ws.Range("K132").Value = 5003.7999
$ws.Range("M132").Value = -2473.7999
$ws.Range("H136").Value = 2975.9583
$ws.Range("I136").Value = 2089.0625
$ws.Range("K136").Value = 6267.1875
$ws.Range("M136").Value = -3717.1875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Value updates
$ws.Range("H94").Value = 90910504
$ws.Range("I94").Value = 105264370
$ws.Range("J94").Value = 2711.3333
$ws.Range("K94").Value = 105264370
$ws.Range("L94").Value = 2711.3333
$ws.Range("M94").Value = -105263919
$ws.Range("N94").Value = -3613.3333
$ws.Range("H134").Value = 2535.9023
$ws.Range("I134").Value = 2173.3572
$ws.Range("K134").Value = 6520.071599999999
$ws.Range("M134").Value = -3985.071599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Value updates
$ws.Range("H31").Value = 3128332.2
$ws.Range("J31").Value = 7817130
$ws.Range("L31").Value = 7817130
$ws.Range("N31").Value = -7817720
$ws.Range("H34").Value = 3128332.2
$ws.Range("J34").Value = 7817130
$ws.Range("L34").Value = 7817130
$ws.Range("N34").Value = -7817534
$ws.Range("H58").Value = 2597.0952
$ws.Range("I58").Value = 1943.4445
$ws.Range("J58").Value = 3087.3333
$ws.Range("K58").Value = 1943.4445
$ws.Range("L58").Value = 3087.3333
$ws.Range("M58").Value = -1740.4445
$ws.Range("N58").Value = -3493.3333
$ws.Range("H62").Value = 100000000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 100000000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H68").Value = 75295
$ws.Range("J68").Value = 75295
$ws.Range("L68").Value = 75295
$ws.Range("H71").Value = 75295
$ws.Range("J71").Value = 75295
$ws.Range("L71").Value = 225885
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H134").Value = 2396.5557
$ws.Range("I134").Value = 2191.4866
$ws.Range("J134").Value = 3345
$ws.Range("K134").Value = 6574.459800000001
$ws.Range("L134").Value = 10035
$ws.Range("M134").Value = -4039.459800000001
$ws.Range("N134").Value = -15105
$ws.Range("H136").Value = 2597.0952
$ws.Range("I136").Value = 1943.4445
$ws.Range("J136").Value = 3087.3333
$ws.Range("K136").Value = 5830.333500000001
$ws.Range("L136").Value = 9261.999899999999
$ws.Range("M136").Value = -3280.333500000001
$ws.Range("N136").Value = -14361.9999

# New cells (previously empty)
$ws.Range("N68").Value = -76793
$ws.Range("N71").Value = -233373

# Cells cleared (now empty)
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Value updates
$ws.Range("H51").Value = 3339
$ws.Range("J51").Value = 5500
$ws.Range("L51").Value = 16500
$ws.Range("N51").Value = -17420

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Value updates
$ws.Range("H107").Value = 3916.4092
$ws.Range("I107").Value = 329.33334
$ws.Range("K107").Value = 329.33334
$ws.Range("M107").Value = 1590.66666
$ws.Range("H132").Value = 2216.8696
$ws.Range("I132").Value = 1809.3334
$ws.Range("K132").Value = 5428.0002
$ws.Range("M132").Value = -2898.0002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Value updates
$ws.Range("H6").Value = 28193.6
$ws.Range("J6").Value = 28193.6
$ws.Range("L6").Value = 28193.6
$ws.Range("N6").Value = -28417.6
$ws.Range("H16").Value = 1395.05
$ws.Range("I16").Value = 744.5
$ws.Range("J16").Value = 7250
$ws.Range("K16").Value = 744.5
$ws.Range("L16").Value = 7250
$ws.Range("M16").Value = -574.5
$ws.Range("N16").Value = -7590
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H122").Value = 5597.552
$ws.Range("I122").Value = 3469.9546
$ws.Range("J122").Value = 12284.286
$ws.Range("K122").Value = 10409.8638
$ws.Range("L122").Value = 36852.858
$ws.Range("M122").Value = -7959.863799999999
$ws.Range("N122").Value = -41752.858
$ws.Range("H132").Value = 3922.7144
$ws.Range("I132").Value = 3698.0557
$ws.Range("K132").Value = 11094.1671
$ws.Range("M132").Value = -8564.167099999999
$ws.Range("H136").Value = 2805.4
$ws.Range("I136").Value = 2568.5
$ws.Range("K136").Value = 7705.5
$ws.Range("M136").Value = -5155.5

# Cells cleared (now empty)
$ws.Range("N118").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Value updates
$ws.Range("H46").Value = 39049.9
$ws.Range("J46").Value = 39049.9
$ws.Range("L46").Value = 39049.9
$ws.Range("N46").Value = -39511.9
$ws.Range("H52").Value = 14920.5
$ws.Range("I52").Value = 14920.5
$ws.Range("K52").Value = 14920.5
$ws.Range("M52").Value = -14694.5
$ws.Range("H122").Value = 9616794
$ws.Range("I122").Value = 1481.3182
$ws.Range("K122").Value = 4443.9546
$ws.Range("M122").Value = -1993.9546
$ws.Range("H132").Value = 3194.2942
$ws.Range("I132").Value = 2706.5
$ws.Range("J132").Value = 10999
$ws.Range("K132").Value = 8119.5
$ws.Range("L132").Value = 32997
$ws.Range("M132").Value = -5589.5
$ws.Range("N132").Value = -38057
$ws.Range("H134").Value = 39049.9
$ws.Range("J134").Value = 39049.9
$ws.Range("L134").Value = 117149.7
$ws.Range("N134").Value = -122219.7

